# Apply updated team-specific time/transition matrix values
# (Mississippi St._B) per the team_specific_matrix data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1952054794520548
$ws.Range("C2").Value = 0.5787671232876712
$ws.Range("J2").Value = 0.01027397260273973
$ws.Range("P2").Value = 0.1301369863013699
$ws.Range("S2").Value = 0.08561643835616438
$ws.Range("B3").Value = 0.005714285714285714
$ws.Range("C3").Value = 0.02285714285714286
$ws.Range("J3").Value = 0.005714285714285714
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1657142857142857
$ws.Range("J4").Value = 0.03571428571428571
$ws.Range("P4").Value = 0.6607142857142857
$ws.Range("S4").Value = 0.3035714285714285
$ws.Range("B6").Value = 0.08144796380090498
$ws.Range("D6").Value = 0.01357466063348416
$ws.Range("F6").Value = 0.04072398190045249
$ws.Range("J6").Value = 0.2262443438914027
$ws.Range("O6").Value = 0.01809954751131222
$ws.Range("Q6").Value = 0.1538461538461539
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.4072398190045249
$ws.Range("B7").Value = 0.1031390134529148
$ws.Range("D7").Value = 0.02690582959641256
$ws.Range("F7").Value = 0.07174887892376682
$ws.Range("J7").Value = 0.09865470852017937
$ws.Range("O7").Value = 0.01345291479820628
$ws.Range("Q7").Value = 0.1569506726457399
$ws.Range("R7").Value = 0.07174887892376682
$ws.Range("S7").Value = 0.4573991031390134
$ws.Range("B8").Value = 0.07472527472527472
$ws.Range("D8").Value = 0.03076923076923077
$ws.Range("F8").Value = 0.05934065934065934
$ws.Range("J8").Value = 0.1142857142857143
$ws.Range("O8").Value = 0.01758241758241758
$ws.Range("Q8").Value = 0.167032967032967
$ws.Range("R8").Value = 0.08791208791208792
$ws.Range("S8").Value = 0.4483516483516484
$ws.Range("B9").Value = 0.1128205128205128
$ws.Range("D9").Value = 0.01538461538461539
$ws.Range("F9").Value = 0.05128205128205128
$ws.Range("J9").Value = 0.1025641025641026
$ws.Range("O9").Value = 0.01538461538461539
$ws.Range("Q9").Value = 0.1384615384615385
$ws.Range("R9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.4871794871794872
$ws.Range("B10").Value = 0.1032
$ws.Range("D10").Value = 0.0256
$ws.Range("F10").Value = 0.07439999999999999
$ws.Range("J10").Value = 0.1176
$ws.Range("O10").Value = 0.016
$ws.Range("Q10").Value = 0.1936
$ws.Range("R10").Value = 0.0824
$ws.Range("S10").Value = 0.3872
$ws.Range("G11").Value = 0.1396648044692737
$ws.Range("J11").Value = 0.1145251396648045
$ws.Range("K11").Value = 0.2374301675977654
$ws.Range("L11").Value = 0.4860335195530726
$ws.Range("S11").Value = 0.0223463687150838
$ws.Range("G12").Value = 0.7458563535911602
$ws.Range("J12").Value = 0.1823204419889503
$ws.Range("K12").Value = 0.005524861878453038
$ws.Range("L12").Value = 0.04419889502762431
$ws.Range("S12").Value = 0.02209944751381215
$ws.Range("G13").Value = 0.5873015873015873
$ws.Range("J13").Value = 0.3650793650793651
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F14").Value = 0.1428571428571428
$ws.Range("G14").Value = 0.7142857142857143
$ws.Range("J14").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.03349282296650718
$ws.Range("H15").Value = 0.09569377990430622
$ws.Range("I15").Value = 0.07655502392344497
$ws.Range("J15").Value = 0.3157894736842105
$ws.Range("K15").Value = 0.1004784688995215
$ws.Range("M15").Value = 0.02870813397129187
$ws.Range("N15").Value = 0.004784688995215311
$ws.Range("O15").Value = 0.06698564593301436
$ws.Range("S15").Value = 0.277511961722488
$ws.Range("F16").Value = 0.009569377990430622
$ws.Range("H16").Value = 0.1578947368421053
$ws.Range("I16").Value = 0.0861244019138756
$ws.Range("J16").Value = 0.4162679425837321
$ws.Range("K16").Value = 0.1004784688995215
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("N16").Value = 0.004784688995215311
$ws.Range("O16").Value = 0.02870813397129187
$ws.Range("S16").Value = 0.1770334928229665
$ws.Range("F17").Value = 0.0121654501216545
$ws.Range("H17").Value = 0.1849148418491484
$ws.Range("I17").Value = 0.1046228710462287
$ws.Range("J17").Value = 0.3795620437956204
$ws.Range("K17").Value = 0.08759124087591241
$ws.Range("M17").Value = 0.024330900243309
$ws.Range("O17").Value = 0.04866180048661801
$ws.Range("S17").Value = 0.1581508515815085
$ws.Range("F18").Value = 0.02185792349726776
$ws.Range("H18").Value = 0.1693989071038251
$ws.Range("I18").Value = 0.1038251366120219
$ws.Range("J18").Value = 0.3715846994535519
$ws.Range("K18").Value = 0.09289617486338798
$ws.Range("M18").Value = 0.0273224043715847
$ws.Range("N18").Value = 0.00546448087431694
$ws.Range("O18").Value = 0.06557377049180328
$ws.Range("S18").Value = 0.1420765027322404
$ws.Range("F19").Value = 0.01345609065155807
$ws.Range("H19").Value = 0.2124645892351275
$ws.Range("I19").Value = 0.07082152974504249
$ws.Range("J19").Value = 0.3505665722379603
$ws.Range("K19").Value = 0.1267705382436261
$ws.Range("M19").Value = 0.028328611898017
$ws.Range("N19").Value = 0.002124645892351275
$ws.Range("O19").Value = 0.06515580736543909
$ws.Range("S19").Value = 0.1303116147308782
